# Auto-generated edit script to apply cryptos.xlsx price/volume/coin updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''25.907.43'
$ws.Range("E2").Value = '''  -0.54%  '
$ws.Range("D3").Value = '''1.640.99'
$ws.Range("E3").Value = '''  -0.11%  '
$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '''  -0.60%  '
$ws.Range("D5").Value = '''215.18'
$ws.Range("E5").Value = '''  -0.16%  '
$ws.Range("D6").Value = '''0.5049'
$ws.Range("E6").Value = '''  -0.05%  '
$ws.Range("D7").Value = '''1.004'
$ws.Range("E7").Value = '''  -0.58%  '
$ws.Range("D8").Value = '''0.2567'
$ws.Range("E8").Value = '''  -0.51%  '
$ws.Range("D9").Value = '''0.06386'
$ws.Range("E9").Value = '''  -0.78%  '
$ws.Range("D10").Value = '''19.61'
$ws.Range("E10").Value = '''  +0.73%  '
$ws.Range("D11").Value = '''0.07791'
$ws.Range("E11").Value = '''  +0.64%  '
$ws.Range("D12").Value = '''1.654.11'
$ws.Range("E12").Value = '''  +0.63%  '
$ws.Range("D13").Value = '''4.277'
$ws.Range("E13").Value = '''  +0.52%  '
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = '''0.5425'
$ws.Range("E14").Value = '''  -0.59%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '''0.0₅7859'
$ws.Range("E15").Value = '''  -0.77%  '
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = '''64.78'
$ws.Range("E16").Value = '''  +1.80%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '''25.958.70'
$ws.Range("E17").Value = '''  -0.35%  '
$ws.Range("B18").Value = 'Dai'
$ws.Range("C18").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D18").Value = '''1.004'
$ws.Range("E18").Value = '''  -0.58%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '''198.30'
$ws.Range("E19").Value = '''  -2.85%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '''4.386'
$ws.Range("E20").Value = '''  +1.91%  '
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D21").Value = '''9.960'
$ws.Range("E21").Value = '''  -0.44%  '
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").Value = '''5.979'
$ws.Range("E22").Value = '''  +0.16%  '
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").Value = '''1.006'
$ws.Range("E23").Value = '''  -0.49%  '
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = '''1.868'
$ws.Range("E24").Value = '''  -3.28%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '''139.97'
$ws.Range("E25").Value = '''  -1.28%  '
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = '''0.1141'
$ws.Range("E26").Value = '''  -1.32%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '''6.847'
$ws.Range("E27").Value = '''  +1.39%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '''15.72'
$ws.Range("E28").Value = '''  -0.37%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = '''1.243'
$ws.Range("E29").Value = '''  +0.01%  '
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = '''0.04922'
$ws.Range("E30").Value = '''  -2.88%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '''3.262'
$ws.Range("E31").Value = '''  +0.31%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '''3.195'
$ws.Range("E32").Value = '''  -0.01%  '
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").Value = '''1.532'
$ws.Range("E33").Value = '''  -0.74%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '''2.370'
$ws.Range("E34").Value = '''  +1.15%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '''0.8932'
$ws.Range("E35").Value = '''  -0.38%  '
$ws.Range("B36").Value = 'MXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D36").Value = '''2.606'
$ws.Range("E36").Value = '''  -0.57%  '
$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '''1.141.16'
$ws.Range("E37").Value = '''  -0.51%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = '''0.5548'
$ws.Range("E38").Value = '''  -1.61%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.01561'
$ws.Range("E39").Value = '''  -0.81%  '
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").Value = '''1.006'
$ws.Range("E40").Value = '''  -0.45%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '''5.687'
$ws.Range("E41").Value = '''  +0.24%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '''0.8210'
$ws.Range("E42").Value = '''  +0.69%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '''99.40'
$ws.Range("E43").Value = '''  -0.54%  '
$ws.Range("B44").Value = 'RocketPoolETH'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").Value = '''1.779.74'
$ws.Range("E44").Value = '''  -0.01%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '''0.0₈119'
$ws.Range("E45").Value = '''  +4.14%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = '''0.4518'
$ws.Range("E46").Value = '''  -0.41%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '''55.30'
$ws.Range("E47").Value = '''  +0.52%  '
$ws.Range("B48").Value = 'Frax'
$ws.Range("C48").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D48").Value = '''1.007'
$ws.Range("E48").Value = '''  -0.46%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '''0.05072'
$ws.Range("E49").Value = '''  +0.59%  '
$ws.Range("B50").Value = 'USDD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D50").Value = '''1.008'
$ws.Range("E50").Value = '''  -0.14%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = '''0.09501'
$ws.Range("E51").Value = '''  +1.80%  '
